$wb = $excel.ActiveWorkbook

# --- binek sheet: fix indirimli_yillik_faiz (B5) 0.1 -> 0.08 ---
$binek = $wb.Worksheets.Item("binek")
$binek.Range("B5").Value = 0.08
$binek.Range("K20").Select()

# --- LCV sheet: fix indirimli_yillik_faiz (B5) 0.1 -> 0.08 ---
$lcv = $wb.Worksheets.Item("LCV")
$lcv.Range("B5").Value = 0.08
$lcv.Activate()
$lcv.Range("A1:B7").Select()

# --- add new HDV sheet (first HDV model), placed after LCV ---
$hdv = $wb.Worksheets.Add($null, $lcv)
$hdv.Name = "HDV"

$hdv.Range("A1").Value = "degisken"
$hdv.Range("B1").Value = "deger"
$hdv.Range("A1:B1").Font.Bold = $true

$hdv.Range("A2").Value = "hdv_kredi_orani"
$hdv.Range("B2").Value = 0.7

$hdv.Range("A3").Value = "hdv_ortalama_vade (yil)"
$hdv.Range("B3").Value = 3

$hdv.Range("A4").Value = "mevcut_yillik faiz"
$hdv.Range("B4").Value = 0.13
$hdv.Range("B4").NumberFormat = "0%"

$hdv.Range("A5").Value = "indirimli_yillik_faiz"
$hdv.Range("B5").Value = 0.08
$hdv.Range("B5").NumberFormat = "0%"

$hdv.Range("A6").Value = "hdv_max indirimli kredi miktari"
$hdv.Range("B6").Value = 120000
$hdv.Range("B6").NumberFormat = '_-* #,##0_-;\-* #,##0_-;_-* "-"??_-;_-@_-'

$hdv.Range("A7").Value = "kredi_kullanan_hdv_orani"
$hdv.Range("B7").Value = 0.52

$hdv.Range("B7").Select()
